$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value2 = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "63.157.15"
Set-TextValue "E2" "  -1.33%  "
Set-TextValue "D3" "3.057.19"
Set-TextValue "E3" "  -3.02%  "
Set-TextValue "D5" "589.54"
Set-TextValue "E5" "  -0.76%  "
Set-TextValue "D6" "153.29"
Set-TextValue "E6" "  +5.26%  "
Set-TextValue "E7" "  -0.16%  "
Set-TextValue "D8" "0.546"
Set-TextValue "E8" "  +3.07%  "
Set-TextValue "D9" "3.065.83"
Set-TextValue "E9" "  -2.35%  "
Set-TextValue "D10" "0.158"
Set-TextValue "E10" "  -2.35%  "
Set-TextValue "D11" "5.84"
Set-TextValue "E11" "  -0.27%  "
Set-TextValue "E12" "  -0.02%  "
Set-TextValue "B13" "Avalanche"
Set-TextValue "C13" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D13" "37.35"
Set-TextValue "E13" "  +0.39%  "
Set-TextValue "B14" "ShibaInu"
Set-TextValue "C14" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000241"
Set-TextValue "E14" "  -2.73%  "
Set-TextValue "E15" "  -1.83%  "
Set-TextValue "D16" "3.563.29"
Set-TextValue "E16" "  -3.03%  "
Set-TextValue "D17" "7.21"
Set-TextValue "E17" "  -1.22%  "
Set-TextValue "D18" "63.240.11"
Set-TextValue "E18" "  -0.96%  "
Set-TextValue "D19" "3.063.81"
Set-TextValue "E19" "  -2.74%  "
Set-TextValue "D20" "476.92"
Set-TextValue "E20" "  +2.11%  "
Set-TextValue "D21" "14.63"
Set-TextValue "E21" "  +1.79%  "
Set-TextValue "D22" "0.717"
Set-TextValue "E22" "  -2.18%  "
Set-TextValue "E23" "  +0.23%  "
Set-TextValue "D24" "2.40"
Set-TextValue "E24" "  +3.34%  "
Set-TextValue "E25" "  -0.45%  "
Set-TextValue "D26" "81.02"
Set-TextValue "E26" "  -0.34%  "
Set-TextValue "D27" "10.04"
Set-TextValue "E27" "  +3.27%  "
Set-TextValue "D28" "0.997"
Set-TextValue "E28" "  -0.38%  "
Set-TextValue "B29" "NEARProtocol"
Set-TextValue "C29" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D29" "7.33"
Set-TextValue "E29" "  -0.88%  "
Set-TextValue "B30" "PancakeSwap"
Set-TextValue "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "2.68"
Set-TextValue "E30" "  -1.45%  "
Set-TextValue "E31" "  -0.20%  "
Set-TextValue "D32" "2.20"
Set-TextValue "E32" "  -2.02%  "
Set-TextValue "E33" "  +2.74%  "
Set-TextValue "D34" "27.21"
Set-TextValue "E34" "  -1.83%  "
Set-TextValue "D35" "0.0₃0847"
Set-TextValue "E35" "  +0.43%  "
Set-TextValue "D36" "1.04"
Set-TextValue "E36" "  -2.02%  "
Set-TextValue "D37" "6.10"
Set-TextValue "E37" "  -1.14%  "
Set-TextValue "E38" "  +2.72%  "
Set-TextValue "D39" "2.22"
Set-TextValue "E39" "  -4.39%  "
Set-TextValue "D40" "9.31"
Set-TextValue "E40" "  +0.56%  "
Set-TextValue "D41" "50.58"
Set-TextValue "E41" "  -1.59%  "
Set-TextValue "D42" "444.00"
Set-TextValue "E42" "  -3.05%  "
Set-TextValue "D43" "0.285"
Set-TextValue "E43" "  -2.77%  "
Set-TextValue "E44" "  -2.70%  "
Set-TextValue "B45" "Kaspa"
Set-TextValue "C45" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D45" "0.112"
Set-TextValue "E45" "  +2.64%  "
Set-TextValue "B46" "Arweave"
Set-TextValue "C46" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D46" "40.04"
Set-TextValue "E46" "  +0.19%  "
Set-TextValue "D47" "2.796.77"
Set-TextValue "E47" "  -4.47%  "
Set-TextValue "D48" "131.91"
Set-TextValue "E48" "  +1.75%  "
Set-TextValue "D49" "0.999"
Set-TextValue "E49" "  +0.04%  "
Set-TextValue "D50" "25.21"
Set-TextValue "E50" "  +3.96%  "
Set-TextValue "E51" "  +0.11%  "
